# Apply TS Pada Paatam corrections: fix heading bold/size formatting to match
# the rest of the document, and merge stray proofErr-split runs back together.

$d = $word.ActiveDocument

function Get-FirstParagraphContaining($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    throw ("paragraph containing " + $needle + " not found")
}

$p = Get-FirstParagraphContaining('(ignore those which are already incorporated')
[void]$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3545F94C" w14:textId="77777777" w:rsidR="00316BB1" w:rsidRPr="00113311" w:rsidRDefault="00316BB1" w:rsidP="00316BB1"><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00113311"><w:rPr><w:b/><w:szCs w:val="24"/></w:rPr><w:t>(ignore those which are already incorporated in your book’s version and date)</w:t></w:r></w:p>')

$p = Get-FirstParagraphContaining('Section, Paragraph')
[void]$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5879ACAE" w14:textId="77777777" w:rsidR="00316BB1" w:rsidRPr="00D40DD6" w:rsidRDefault="00316BB1" w:rsidP="00A92D09"><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="00D40DD6"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Section, Paragraph</w:t></w:r></w:p>')

$p = Get-FirstParagraphContaining('Reference')
[void]$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4EE26168" w14:textId="77777777" w:rsidR="00316BB1" w:rsidRPr="00D40DD6" w:rsidRDefault="00316BB1" w:rsidP="00A92D09"><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="00D40DD6"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Reference</w:t></w:r></w:p>')

$p = Get-FirstParagraphContaining('As Printed')
[void]$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5DBFB696" w14:textId="77777777" w:rsidR="00316BB1" w:rsidRPr="00D40DD6" w:rsidRDefault="00316BB1" w:rsidP="00A92D09"><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="00D40DD6"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>As Printed</w:t></w:r></w:p>')

$p = Get-FirstParagraphContaining('To be read as or corrected as')
[void]$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1ECC703C" w14:textId="77777777" w:rsidR="00316BB1" w:rsidRPr="00D40DD6" w:rsidRDefault="00316BB1" w:rsidP="00A92D09"><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:right="-18"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="00D40DD6"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>To be read as or corrected as</w:t></w:r></w:p>')

$p = Get-FirstParagraphContaining('1.4.29.1')
[void]$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0F0EF36C" w14:textId="77777777" w:rsidR="009C7EFC" w:rsidRPr="00380E2C" w:rsidRDefault="00F74342" w:rsidP="00FD1B8A"><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="-18"/><w:rPr><w:rFonts w:cs="Latha"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="ta-IN"/></w:rPr></w:pPr><w:r w:rsidRPr="00380E2C"><w:rPr><w:rFonts w:cs="Latha"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="ta-IN"/></w:rPr><w:t>1.4.29.1  Padam 6</w:t></w:r></w:p>')

$p = Get-FirstParagraphContaining('1.4.1.2  Vaakyam')
[void]$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1748BADC" w14:textId="77777777" w:rsidR="003827E2" w:rsidRPr="00380E2C" w:rsidRDefault="003827E2" w:rsidP="003827E2"><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:left="-108" w:right="-166"/><w:rPr><w:rFonts w:cs="Latha"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="ta-IN"/></w:rPr></w:pPr><w:r w:rsidRPr="00380E2C"><w:rPr><w:rFonts w:cs="Latha"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="ta-IN"/></w:rPr><w:lastRenderedPageBreak/><w:t>TS 1.4.1.2  Vaakyam</w:t></w:r></w:p>')

Write-Output "TS Pada Padam 1.1 to 1.8 final corrections applied."
